# Add new CAN messages
$wb = $excel.ActiveWorkbook

#
# "CAN IDs" sheet: document the new 130-134 input-readings mailbox range
# that sits just above the existing 135-140 raw-ADC-readings block.
$wsIds = $wb.Worksheets.Item("CAN IDs")
$wsIds.Range("F7").Value = "130-134"
$wsIds.Range("H7").Value = "Input Readings"

# Small scratch table (mailbox # / label) used while drafting the new
# analog-input mailboxes below.
$wsIds.Range("P16").Value = 4
$wsIds.Range("Q16").Value = "BP1"
$wsIds.Range("P17").Value = 4
$wsIds.Range("Q17").Value = "BP2"
$wsIds.Range("P19").Value = 4
$wsIds.Range("Q19").Value = "IO1"
$wsIds.Range("P20").Value = 4
$wsIds.Range("Q20").Value = "IO2"
$wsIds.Range("P21").Value = 5
$wsIds.Range("Q21").Value = "IO3"
$wsIds.Range("P22").Value = 5
$wsIds.Range("Q22").Value = "IO4"
$wsIds.Range("P23").Value = 5
$wsIds.Range("Q23").Value = "IO5"
$wsIds.Range("P24").Value = 5
$wsIds.Range("Q24").Value = "IO6"

# "Mailboxes (CANB)" sheet: rename the existing HO_CAN mailbox targets to
# go through the new CAN2 struct, and add the new mailboxes (8-13) that
# now carry the Output/CPU1 analog + mode data.
$wsCanB = $wb.Worksheets.Item("Mailboxes (CANB)")

$wsCanB.Range("E3").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E4").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E5").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E6").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E7").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E8").Value = "HO_CAN->CAN2->rawInputs"
$wsCanB.Range("E9").Value = "HO_CAN->CAN2->modes"
$wsCanB.Range("E10").Value = "HO_CAN->CAN2->motorControlSlow"

$wsCanB.Range("C11").Value = "Output"
$wsCanB.Range("D11").Value = "CPU1"
$wsCanB.Range("E11").Value = "HO_CAN->CAN2->modes"

$wsCanB.Range("C12").Value = "Output"
$wsCanB.Range("D12").Value = "CPU1"
$wsCanB.Range("E12").Value = "HO_CAN->CAN2->analogInputs"

$wsCanB.Range("C13").Value = "Output"
$wsCanB.Range("D13").Value = "CPU1"
$wsCanB.Range("E13").Value = "HO_CAN->CAN2->analogInputs"

$wsCanB.Range("C14").Value = "Output"
$wsCanB.Range("D14").Value = "CPU1"
$wsCanB.Range("E14").Value = "HO_CAN->CAN2->analogInputs"

$wsCanB.Range("C15").Value = "Output"
$wsCanB.Range("D15").Value = "CPU1"
$wsCanB.Range("E15").Value = "HO_CAN->CAN2->analogInputs"

$wsCanB.Range("C16").Value = "Output"
$wsCanB.Range("D16").Value = "CPU1"
$wsCanB.Range("E16").Value = "HO_CAN->CAN2->analogInputs"

# Restore the selection on each affected sheet to match where the author
# ended up after making the edits.
$wsIds.Range("U23").Select()
$wsCanB.Range("J14").Select()
